# fix: clean up sheet
#
# - P&L had a stray, always-empty column B (the month headers/values actually
#   started in column C). Delete the empty column so data starts at B.
# - Give "Salaries by Cost Center" a proper header for its row-label column
#   ("Cost Cetner") and bold the row-label column on both detail sheets so
#   they read consistently with the rest of the workbook.

$wb = $excel.ActiveWorkbook

$wsPL       = $wb.Worksheets.Item("P&L")
$wsSales    = $wb.Worksheets.Item("Sales by Product")
$wsSalaries = $wb.Worksheets.Item("Salaries by Cost Center")

# --- P&L: remove the empty leading data column ------------------------------
$wsPL.Columns("B").Delete()

# --- Salaries by Cost Center: label the row-header column ------------------
$wsSalaries.Range("A1").Value = "Cost Cetner"

# --- Bold the row-label column on both detail sheets ------------------------
$wsSales.Range("A2:A4").Font.Bold = $true
$wsSales.Range("A2:A4").Font.Name = "Aptos Narrow"
$wsSales.Range("A2:A4").Select()

$wsSalaries.Range("A1:A4").Font.Bold = $true
$wsSalaries.Range("A1:A4").Font.Name = "Aptos Narrow"
$wsSalaries.Range("A2:A4").Select()

# --- Leave P&L as the active sheet/selection --------------------------------
$wsPL.Activate()
$wsPL.Range("B3").Select()
